$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 2.395846366882324
$ws.Range("B2").Value = 5.165933609008789
$ws.Range("C2").Value = -3.433014392852783
$ws.Range("A3").Value = -4.131561756134033
$ws.Range("B3").Value = 1.293697118759155
$ws.Range("C3").Value = 3.305315971374512
$ws.Range("A4").Value = -11.94900512695312
$ws.Range("B4").Value = -0.5390903949737549
$ws.Range("C4").Value = 0.9473530054092408
$ws.Range("A5").Value = 6.543986320495605
$ws.Range("B5").Value = -9.976800918579102
$ws.Range("C5").Value = -2.196774959564209
$ws.Range("A6").Value = 2.871486902236938
$ws.Range("B6").Value = -5.876065254211426
$ws.Range("C6").Value = -3.124087810516357
$ws.Range("A7").Value = 7.37382698059082
$ws.Range("B7").Value = 11.17136192321777
$ws.Range("C7").Value = -1.706753373146057
$ws.Range("A8").Value = -4.830907821655273
$ws.Range("B8").Value = 0.3003380000591278
$ws.Range("C8").Value = 1.922602653503418
$ws.Range("A9").Value = -6.295113563537598
$ws.Range("B9").Value = 2.037784337997437
$ws.Range("C9").Value = 3.717573165893554
$ws.Range("A10").Value = -1.753358721733093
$ws.Range("B10").Value = 12.93171119689941
$ws.Range("C10").Value = 2.580935955047607
$ws.Range("A11").Value = -3.18294358253479
$ws.Range("B11").Value = -2.735665798187256
$ws.Range("C11").Value = 0.6512094736099243
$ws.Range("A12").Value = 7.503256797790527
$ws.Range("B12").Value = 0.8255133628845215
$ws.Range("C12").Value = -1.698231339454651
$ws.Range("A13").Value = -1.000749349594116
$ws.Range("B13").Value = -2.794787883758545
$ws.Range("C13").Value = -1.550159573554993
$ws.Range("A14").Value = -2.628473520278931
$ws.Range("B14").Value = -2.876280546188354
$ws.Range("C14").Value = -0.96160089969635
$ws.Range("A15").Value = -5.674064636230469
$ws.Range("B15").Value = -1.664009690284729
$ws.Range("C15").Value = -3.089466571807861
$ws.Range("A16").Value = 0.2152500003576278
$ws.Range("B16").Value = 15.8856897354126
$ws.Range("C16").Value = -3.721701145172119
$ws.Range("A17").Value = 1.613942265510559
$ws.Range("B17").Value = -4.874184131622314
$ws.Range("C17").Value = -0.0044607948511838
$ws.Range("A18").Value = 4.66752290725708
$ws.Range("B18").Value = -2.795853137969971
$ws.Range("C18").Value = -8.249074935913086
$ws.Range("A19").Value = 1.00194776058197
$ws.Range("B19").Value = -2.820886850357056
$ws.Range("C19").Value = 4.437692165374756
$ws.Range("A20").Value = -8.889565467834473
$ws.Range("B20").Value = 1.164800047874451
$ws.Range("C20").Value = 11.45938301086426
$ws.Range("A21").Value = 2.723947763442993
$ws.Range("B21").Value = 6.857040882110596
$ws.Range("C21").Value = -5.117730140686035
$ws.Range("A22").Value = 1.208609104156494
$ws.Range("B22").Value = -8.24947452545166
$ws.Range("C22").Value = -0.5370930433273315
$ws.Range("A23").Value = 8.961604118347168
$ws.Range("B23").Value = 6.756906032562256
$ws.Range("C23").Value = -2.651110410690308
$ws.Range("A24").Value = -2.561361789703369
$ws.Range("B24").Value = 1.376787781715393
$ws.Range("C24").Value = 1.834718346595764
$ws.Range("A25").Value = -2.701976776123047
$ws.Range("B25").Value = 1.867341995239258
$ws.Range("C25").Value = 1.237104892730713
$ws.Range("A26").Value = -5.291634559631348
$ws.Range("B26").Value = -3.584681510925293
$ws.Range("C26").Value = -2.499310255050659
$ws.Range("A27").Value = -0.4974119365215301
$ws.Range("B27").Value = 12.873122215271
$ws.Range("C27").Value = 5.119993686676025
$ws.Range("A28").Value = -0.0159123875200748
$ws.Range("B28").Value = 1.14615797996521
$ws.Range("C28").Value = 3.750063896179199
$ws.Range("A29").Value = 8.204200744628906
$ws.Range("B29").Value = 6.466621398925781
$ws.Range("C29").Value = -1.093693733215332
$ws.Range("A30").Value = -1.878527283668518
$ws.Range("B30").Value = -2.402770519256592
$ws.Range("C30").Value = 0.963331937789917
$ws.Range("A31").Value = -2.405300617218018
$ws.Range("B31").Value = 1.181844353675843
$ws.Range("C31").Value = 4.602807998657227
